$wb = $excel.ActiveWorkbook

# --- CPPbI sheet: update CO2 capture-ratio values in columns B & C (rows 2-26) ---
$ws = $wb.Worksheets.Item("CPPbI")

$newValues = @(
    0.9,  # row 2  - agriculture and forestry 01T03
    0.9,  # row 3  - coal mining 05
    0.9,  # row 4  - oil and gas extraction 06
    0.9,  # row 5  - other mining and quarrying 07T08
    0,    # row 6  - food beverage and tobacco 10T12
    0,    # row 7  - textiles apparel and leather 13T15
    0,    # row 8  - wood products 16
    0.9,  # row 9  - pulp paper and printing 17T18
    0.9,  # row 10 - refined petroleum and coke 19
    0.9,  # row 11 - chemicals 20
    0,    # row 12 - rubber and plastic products 22
    0.9,  # row 13 - glass and glass products 231
    0.9,  # row 14 - cement and other nonmetallic minerals 239
    0.9,  # row 15 - iron and steel 241
    0,    # row 16 - other metals 242
    0,    # row 17 - metal products except machinery and vehicles 25
    0,    # row 18 - computers and electronics 26
    0,    # row 19 - appliances and electrical equipment 27
    0,    # row 20 - other machinery 28
    0,    # row 21 - road vehicles 29
    0.9,  # row 22 - nonroad vehicles 30
    0,    # row 23 - other manufacturing 31T33
    0.9,  # row 24 - energy pipelines and gas processing 352T353
    0.9,  # row 25 - water and waste 36T39
    0.9   # row 26 - construction 41T43
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
    $ws.Cells.Item($row, 3).Value = $newValues[$i]
}

# Widen column B so the new values are readable
$ws.Columns.Item(2).ColumnWidth = 28

# Make CPPbI the active/visible sheet and set its selection
$ws.Activate() | Out-Null
$ws.Range("C18").Select() | Out-Null
